function Set-DatePlaceholderText {
    param($shapes, $newText)
    $count = $shapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
        }
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$p = $ppt.ActivePresentation
$newDate = "1/12/2021"

# Slide master date placeholder
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Every slide layout's date placeholder
$layouts = $p.SlideMaster.CustomLayouts
$layoutCount = $layouts.Count
for ($li = 1; $li -le $layoutCount; $li++) {
    $lay = $layouts.Item($li)
    Set-DatePlaceholderText $lay.Shapes $newDate
}

# Notes master date placeholder
Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# Resize "Rectangle 4" on the last slide (slide 16)
$lastSlide = $p.Slides.Item($p.Slides.Count)
$shapes = $lastSlide.Shapes
$shapeCount = $shapes.Count
for ($i = 1; $i -le $shapeCount; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.Name -eq "Rectangle 4") {
        $sh.Width = 1833563
    }
}
